$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C (Telefone) and D (idPagamento) hold digit-only strings that
# must stay text (not be coerced into numbers). Leading "'" forces Excel
# to store the value as text - for D10:D13 there is no idPagamento yet,
# so after the "'" strips off we are left with an empty text string,
# matching the source rows' "" placeholder.

# --- Row 10 ---
$ws.Range("A10").Value = "Vitor Ito"
$ws.Range("B10").Value = 1578424633
$ws.Range("C10").Value = "'11966548087"
$ws.Range("D10").Value = "'"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 3
$ws.Range("H10").Value = 4
$ws.Range("I10").Value = 5
$ws.Range("J10").Value = 6
$ws.Range("K10").Value = 9
$ws.Range("L10").Value = 10
$ws.Range("M10").Value = 11
$ws.Range("N10").Value = 38
$ws.Range("O10").Value = "Não"

# --- Row 11 ---
$ws.Range("A11").Value = "Vitor Ito"
$ws.Range("B11").Value = 1578424633
$ws.Range("C11").Value = "'11900009999"
$ws.Range("D11").Value = "'"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 3
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = 5
$ws.Range("J11").Value = 6
$ws.Range("K11").Value = 7
$ws.Range("L11").Value = 8
$ws.Range("M11").Value = 9
$ws.Range("N11").Value = 10
$ws.Range("O11").Value = "Não"

# --- Row 12 ---
$ws.Range("A12").Value = "Vitor Ito"
$ws.Range("B12").Value = 1578424633
$ws.Range("C12").Value = "'11966548087"
$ws.Range("D12").Value = "'"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 4
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 6
$ws.Range("K12").Value = 7
$ws.Range("L12").Value = 8
$ws.Range("M12").Value = 9
$ws.Range("N12").Value = 10
$ws.Range("O12").Value = "Não"

# --- Row 13 ---
$ws.Range("A13").Value = "Vitor Ito"
$ws.Range("B13").Value = 1578424633
$ws.Range("C13").Value = "'11966548087"
$ws.Range("D13").Value = "'"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 3
$ws.Range("H13").Value = 4
$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 6
$ws.Range("K13").Value = 7
$ws.Range("L13").Value = 8
$ws.Range("M13").Value = 9
$ws.Range("N13").Value = 10
$ws.Range("O13").Value = "Não"
